$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.00" or
# "0.0000271" are not coerced into numbers and keep their exact formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '98.529.04'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '3.348.18'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '260.86'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").Value = '647.10'
$ws.Range("E6").Value = '  +1.83%  '
$ws.Range("D7").Value = '1.53'
$ws.Range("E7").Value = '  +10.33%  '
$ws.Range("D8").Value = '0.460'
$ws.Range("E8").Value = '  +16.43%  '
$ws.Range("D9").Value = '1.08'
$ws.Range("E9").Value = '  +22.53%  '
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").Value = '3.346.53'
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").Value = '44.09'
$ws.Range("E12").Value = '  +21.33%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.208'
$ws.Range("E13").Value = '  +3.93%  '
$ws.Range("D14").Value = '0.0000271'
$ws.Range("E14").Value = '  +8.53%  '
$ws.Range("D15").Value = '98.558.28'
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").Value = '3.978.09'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").Value = '5.56'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '3.349.75'
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").Value = '7.49'
$ws.Range("E19").Value = '  +20.21%  '
$ws.Range("D20").Value = '16.80'
$ws.Range("E20").Value = '  +10.32%  '
$ws.Range("D21").Value = '539.08'
$ws.Range("E21").Value = '  +8.77%  '
$ws.Range("D22").Value = '3.60'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = '10.17'
$ws.Range("E23").Value = '  +8.00%  '
$ws.Range("D24").Value = '0.0000214'
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("D25").Value = '0.426'
$ws.Range("E25").Value = '  +50.18%  '
$ws.Range("D26").Value = '104.03'
$ws.Range("E26").Value = '  +15.04%  '
$ws.Range("D27").Value = '6.21'
$ws.Range("E27").Value = '  +7.12%  '
$ws.Range("D28").Value = '12.72'
$ws.Range("E28").Value = '  +4.89%  '
$ws.Range("D29").Value = '3.528.25'
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").Value = '0.149'
$ws.Range("E30").Value = '  +10.55%  '
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").Value = '10.91'
$ws.Range("E32").Value = '  +12.53%  '
$ws.Range("D33").Value = '0.187'
$ws.Range("E33").Value = '  -5.91%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").Value = '29.12'
$ws.Range("E35").Value = '  +4.54%  '
$ws.Range("D36").Value = '0.519'
$ws.Range("E36").Value = '  +9.95%  '
$ws.Range("E37").Value = '  +5.56%  '
$ws.Range("D38").Value = '0.157'
$ws.Range("E38").Value = '  +3.80%  '
$ws.Range("E39").Value = '  +3.04%  '
$ws.Range("D40").Value = '516.68'
$ws.Range("E40").Value = '  +1.93%  '
$ws.Range("D41").Value = '24.71'
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("E42").Value = '  +2.87%  '
$ws.Range("D43").Value = '3.83'
$ws.Range("E43").Value = '  +2.72%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '3.35'
$ws.Range("E44").Value = '  -1.51%  '
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").Value = '0.808'
$ws.Range("E45").Value = '  +2.70%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '0.0396'
$ws.Range("E47").Value = '  +21.71%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '2.02'
$ws.Range("E48").Value = '  +3.34%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '163.96'
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("D50").Value = '7.70'
$ws.Range("E50").Value = '  +16.26%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").Value = '4.95'
$ws.Range("E51").Value = '  +5.03%  '
